$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.996.60'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.647.30'
$ws.Range("E3").Value = '  -0.90%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.92'
$ws.Range("E5").Value = '  +2.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5225'
$ws.Range("E6").Value = '  +0.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.39%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2613'
$ws.Range("E8").Value = '  +1.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06350'
$ws.Range("E9").Value = '  +0.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.77'
$ws.Range("E10").Value = '  -0.94%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07683'
$ws.Range("E11").Value = '  +1.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.648.15'
$ws.Range("E12").Value = '  -0.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.429'
$ws.Range("E13").Value = '  +0.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.870.83'
$ws.Range("E14").Value = '  -0.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5507'
$ws.Range("E15").Value = '  +2.76%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8235'
$ws.Range("E16").Value = '  +3.92%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.70'
$ws.Range("E17").Value = '  -2.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.980.44'
$ws.Range("E18").Value = '  -0.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.719'
$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.61'
$ws.Range("E21").Value = '  +1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.21'
$ws.Range("E22").Value = '  +0.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.282'
$ws.Range("E23").Value = '  +1.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.23'
$ws.Range("E25").Value = '  -3.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1246'
$ws.Range("E26").Value = '  +2.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.394'
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.98'
$ws.Range("E28").Value = '  +2.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.410'
$ws.Range("E29").Value = '  +2.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05932'
$ws.Range("E30").Value = '  -3.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.258'
$ws.Range("E31").Value = '  -0.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.424'
$ws.Range("E32").Value = '  -1.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.407'
$ws.Range("E33").Value = '  +0.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.642'
$ws.Range("E34").Value = '  +0.79%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9933'
$ws.Range("E35").Value = '  +0.71%  '

$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.743'
$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5635'
$ws.Range("E38").Value = '  -3.95%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01604'
$ws.Range("E39").Value = '  +0.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.867'
$ws.Range("E40").Value = '  -1.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8562'
$ws.Range("E41").Value = '  +1.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.026.15'
$ws.Range("E43").Value = '  -7.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.93'
$ws.Range("E44").Value = '  -0.89%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.795.47'
$ws.Range("E45").Value = '  -1.01%  '

$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.70'
$ws.Range("E47").Value = '  +1.46%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  +0.28%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.055'
$ws.Range("E49").Value = '  +1.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05140'
$ws.Range("E50").Value = '  -2.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4206'
$ws.Range("E51").Value = '  -0.95%  '
